# Remove the "1999-2001" Saitama prefectural institute row (row 20).
# This shifts rows 21-23 up by one, which naturally produces the new
# contents for rows 20-22 described in the diff, and shrinks the used
# range from A1:C23 to A1:C22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(20).Delete()
